$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 299.8387
$ws.Range("J33").Value = 756.6667
$ws.Range("L33").Value = 756.6667
$ws.Range("N33").Value = -1214.6667
$ws.Range("H38").Value = 1613.3
$ws.Range("I38").Value = 125.888885
$ws.Range("K38").Value = 377.666655
$ws.Range("M38").Value = -5.666654999999992
$ws.Range("H40").Value = 2455.1177
$ws.Range("I40").Value = 2474.625
$ws.Range("J40").Value = 2437.7778
$ws.Range("K40").Value = 2474.625
$ws.Range("L40").Value = 2437.7778
$ws.Range("M40").Value = -2299.625
$ws.Range("N40").Value = -2787.7778
$ws.Range("H70").Value = 1427.4667
$ws.Range("J70").Value = 1421.9166
$ws.Range("L70").Value = 4265.7498
$ws.Range("N70").Value = -4805.7498
$ws.Range("H73").Value = 1427.4667
$ws.Range("J73").Value = 1421.9166
$ws.Range("L73").Value = 4265.7498
$ws.Range("N73").Value = -6137.7498
$ws.Range("H82").Value = 2013.3334
$ws.Range("I82").Value = 2013.3334
$ws.Range("K82").Value = 6040.0002
$ws.Range("M82").Value = -5634.0002
$ws.Range("H85").Value = 2013.3334
$ws.Range("I85").Value = 2013.3334
$ws.Range("K85").Value = 6040.0002
$ws.Range("M85").Value = -4636.0002
$ws.Range("H96").Value = 62501816
$ws.Range("I96").Value = 1860.8
$ws.Range("K96").Value = 5582.4
$ws.Range("M96").Value = -4209.4
$ws.Range("H99").Value = 58829308
$ws.Range("I99").Value = 359
$ws.Range("K99").Value = 1077
$ws.Range("M99").Value = 421
$ws.Range("H100").Value = 10898.846
$ws.Range("I100").Value = 4848
$ws.Range("J100").Value = 16085.286
$ws.Range("K100").Value = 4848
$ws.Range("L100").Value = 16085.286
$ws.Range("M100").Value = -4307
$ws.Range("N100").Value = -17167.286
$ws.Range("H112").Value = 6433.8184
$ws.Range("J112").Value = 5443.8887
$ws.Range("L112").Value = 16331.6661
$ws.Range("N112").Value = -18547.6661
$ws.Range("H127").Value = 3738.5625
$ws.Range("I127").Value = 4129.143
$ws.Range("J127").Value = 1004.5
$ws.Range("K127").Value = 12387.429
$ws.Range("L127").Value = 3013.5
$ws.Range("M127").Value = -7427.429
$ws.Range("N127").Value = -12933.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 603.22
$ws.Range("I32").Value = 471.70526
$ws.Range("J32").Value = 3102
$ws.Range("K32").Value = 471.70526
$ws.Range("L32").Value = 3102
$ws.Range("M32").Value = -184.70526
$ws.Range("N32").Value = -3676
$ws.Range("H61").Value = 7766.6313
$ws.Range("I61").Value = 8970.799999999999
$ws.Range("K61").Value = 8970.799999999999
$ws.Range("M61").Value = -8758.799999999999
$ws.Range("H63").Value = 3311.75
$ws.Range("I63").Value = 1699
$ws.Range("J63").Value = 8150
$ws.Range("K63").Value = 1699
$ws.Range("L63").Value = 8150
$ws.Range("M63").Value = -1013
$ws.Range("N63").Value = -9522
$ws.Range("H66").Value = 3311.75
$ws.Range("I66").Value = 1699
$ws.Range("J66").Value = 8150
$ws.Range("K66").Value = 8495
$ws.Range("L66").Value = 40750
$ws.Range("M66").Value = -5063
$ws.Range("N66").Value = -47614
$ws.Range("H136").Value = 7766.6313
$ws.Range("I136").Value = 8970.799999999999
$ws.Range("K136").Value = 26912.4
$ws.Range("M136").Value = -24362.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H107").Value = 50500000
$ws.Range("I107").Value = 50500000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 50500000
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -50498080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 198.6
$ws.Range("I7").Value = 175.3077
$ws.Range("J7").Value = 350
$ws.Range("K7").Value = 175.3077
$ws.Range("L7").Value = 350
$ws.Range("M7").Value = -62.30770000000001
$ws.Range("N7").Value = -576
$ws.Range("H62").Value = 9053.200000000001
$ws.Range("I62").Value = 5402.5
$ws.Range("J62").Value = 9965.875
$ws.Range("K62").Value = 5402.5
$ws.Range("L62").Value = 9965.875
$ws.Range("M62").Value = -4778.5
$ws.Range("N62").Value = -11213.875
$ws.Range("H65").Value = 9053.200000000001
$ws.Range("I65").Value = 5402.5
$ws.Range("J65").Value = 9965.875
$ws.Range("K65").Value = 27012.5
$ws.Range("L65").Value = 49829.375
$ws.Range("M65").Value = -23892.5
$ws.Range("N65").Value = -56069.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 3300
$ws.Range("I17").Value = 3000
$ws.Range("K17").Value = 9000
$ws.Range("M17").Value = -8831
$ws.Range("H56").Value = 5207.909
$ws.Range("I56").Value = 5207.909
$ws.Range("K56").Value = 5207.909
$ws.Range("M56").Value = -4677.909
$ws.Range("H70").Value = 5747.125
$ws.Range("I70").Value = 5747.125
$ws.Range("K70").Value = 17241.375
$ws.Range("M70").Value = -16926.375
$ws.Range("H73").Value = 5747.125
$ws.Range("I73").Value = 5747.125
$ws.Range("K73").Value = 17241.375
$ws.Range("M73").Value = -16149.375
$ws.Range("H95").Value = 12000
$ws.Range("I95").Value = 4500
$ws.Range("J95").Value = 19500
$ws.Range("K95").Value = 13500
$ws.Range("L95").Value = 58500
$ws.Range("M95").Value = -11441
$ws.Range("N95").Value = -62618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11816.333
$ws.Range("I122").Value = 11516.333
$ws.Range("J122").Value = 12416.333
$ws.Range("K122").Value = 34548.999
$ws.Range("L122").Value = 37248.999
$ws.Range("M122").Value = -32098.999
$ws.Range("N122").Value = -42148.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 3013.9
$ws.Range("I55").Value = 915.2
$ws.Range("J55").Value = 5112.6
$ws.Range("K55").Value = 915.2
$ws.Range("L55").Value = 5112.6
$ws.Range("M55").Value = -742.2
$ws.Range("N55").Value = -5458.6
$ws.Range("H68").Value = 1013.0833
$ws.Range("I68").Value = 808
$ws.Range("J68").Value = 1423.25
$ws.Range("K68").Value = 808
$ws.Range("L68").Value = 1423.25
$ws.Range("M68").Value = -59
$ws.Range("N68").Value = -2921.25
$ws.Range("H71").Value = 1013.0833
$ws.Range("I71").Value = 808
$ws.Range("J71").Value = 1423.25
$ws.Range("K71").Value = 4040
$ws.Range("L71").Value = 7116.25
$ws.Range("M71").Value = -296
$ws.Range("N71").Value = -14604.25
$ws.Range("H137").Value = 150000
$ws.Range("J137").Value = 150000
$ws.Range("L137").Value = 150000
$ws.Range("M137").Value = -160200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9982
$ws.Range("J62").Value = 11620.4
$ws.Range("L62").Value = 11620.4
$ws.Range("N62").Value = -12868.4
$ws.Range("H65").Value = 9982
$ws.Range("J65").Value = 11620.4
$ws.Range("L65").Value = 58102
$ws.Range("N65").Value = -64342
